$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the weekly figures between row 2 (week of 2023-02-16) and row 3
# (week of 2023-02-23): dates, volume, prices and $/Kg.

$ws.Range("D2").Value = 44980
$ws.Range("M2").Value = 50
$ws.Range("N2").Value = 25000
$ws.Range("O2").Value = 25000
$ws.Range("P2").Value = 25000
$ws.Range("S2").Value = 3125

$ws.Range("D3").Value = 44973
$ws.Range("M3").Value = 55
$ws.Range("N3").Value = 28000
$ws.Range("O3").Value = 28000
$ws.Range("P3").Value = 28000
$ws.Range("S3").Value = 3500
